$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 17.10987091064453
$ws.Range("D2").Value = 181

$ws.Range("C3").Value = 16.53409004211426
$ws.Range("D3").Value = 175

$ws.Range("C4").Value = 16.64018630981445
$ws.Range("D4").Value = 196

$ws.Range("C5").Value = 16.49308204650879
$ws.Range("D5").Value = 184

$ws.Range("C6").Value = 16.56007766723633
$ws.Range("D6").Value = 184
